$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.717.45"
$ws.Range("E2").Value = "  +3.10%  "
$ws.Range("D3").Value = "3.973.11"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'585.55"
$ws.Range("E5").Value = "  +8.90%  "
$ws.Range("D6").Value = "'157.08"
$ws.Range("E6").Value = "  +6.74%  "
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").Value = "'53.30"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "'0.0000319"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "'10.79"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").Value = "4.610.48"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "3.966.99"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("E16").Value = "  +9.37%  "
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "'20.36"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D20").Value = "72.388.72"
$ws.Range("E20").Value = "  +2.86%  "
$ws.Range("D21").Value = "'432.07"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "'4.69"
$ws.Range("E22").Value = "  +12.31%  "
$ws.Range("D23").Value = "'95.92"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").Value = "'14.33"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "'4.45"
$ws.Range("E26").Value = "  +22.80%  "
$ws.Range("D27").Value = "'11.07"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "'10.65"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("D30").Value = "'36.45"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "'7.78"
$ws.Range("E31").Value = "  +4.71%  "
$ws.Range("D32").Value = "'13.56"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("E33").Value = "  +2.12%  "
$ws.Range("D34").Value = "'677.70"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "'48.42"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("D36").Value = "'68.81"
$ws.Range("E36").Value = "  +6.07%  "
$ws.Range("D37").Value = "0.0₃0875"
$ws.Range("E37").Value = "  +7.14%  "
$ws.Range("D38").Value = "'0.435"
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.146"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'3.33"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'0.0486"
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("D45").Value = "'10.74"
$ws.Range("E45").Value = "  +12.69%  "
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("D47").Value = "'2.64"
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("D48").Value = "'3.38"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").Value = "'3.02"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").Value = "'3.41"
$ws.Range("E50").Value = "  +5.64%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.796.43"
$ws.Range("E51").Value = "  +12.40%  "
